{"js": "const body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\nlet items = body.paragraphs.items;\n\n// 1. Fix typo: \"wither\" -> \"either\" in the \"Unable to figure out...\" paragraph.\nconst imgIssuePara = items[15];\nconst witherResults = imgIssuePara.search(\"wither\", { matchCase: true });\nwitherResults.load(\"items\");\nawait context.sync();\nif (witherResults.items.length > 0) {\n  witherResults.items[0].insertText(\"either\", \"Replace\");\n  await context.sync();\n}\n\n// 2. Remove the blank paragraph that separated the image paragraph from the\n//    \"Known Bugs\" heading, effectively moving \"Known Bugs\" directly under it.\nitems[16].delete();\nawait context.sync();\n\n// Paragraph indices shifted by -1 from here on because of the deletion above.\nbody.paragraphs.load(\"items\");\nawait context.sync();\nitems = body.paragraphs.items;\n\n// 3. \"Sometimes the keystore will be wiped\" -> \"Sometimes the keystore file will be wiped\"\nconst wipeKeystorePara = items[19];\nconst keystoreResults = wipeKeystorePara.search(\"keystore will be wiped\", { matchCase: true });\nkeystoreResults.load(\"items\");\nawait context.sync();\nif (keystoreResults.items.length > 0) {\n  keystoreResults.items[0].insertText(\"keystore file will be wiped\", \"Replace\");\n  await context.sync();\n}\n\n// 4. \"save the file\" -> \"save the updated keystore file\"\nconst saveFilePara = items[20];\nconst saveFileResults = saveFilePara.search(\"save the file\", { matchCase: true });\nsaveFileResults.load(\"items\");\nawait context.sync();\nif (saveFileResults.items.length > 0) {\n  saveFileResults.items[0].insertText(\"save the updated keystore file\", \"Replace\");\n  await context.sync();\n}\n\n// 5. \"temperature the humidity\" -> \"temperature, the humidity\"\nconst humidityPara = items[24];\nconst commaResults = humidityPara.search(\"temperature the humidity\", { matchCase: true });\ncommaResults.load(\"items\");\nawait context.sync();\nif (commaResults.items.length > 0) {\n  commaResults.items[0].insertText(\"temperature, the humidity\", \"Replace\");\n  await context.sync();\n}\n\n// 6. \"so it also outputs\" -> \"so it will output\"\nconst outputResults = humidityPara.search(\"so it also outputs\", { matchCase: true });\noutputResults.load(\"items\");\nawait context.sync();\nif (outputResults.items.length > 0) {\n  outputResults.items[0].insertText(\"so it will output\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Fix typo: \"wither\" -> \"either\" in the \"Unable to figure out...\" paragraph\n#    (paragraph right after the \"Image sending and retrieval.\" heading).\n$witherPara = $d.Paragraphs.Item(16)\n[void]$witherPara.Range.Find.Execute(\"wither\", $false, $false, $false, $false, $false, $true, 1, $false, \"either\", 2)\n\n# 2. Remove the blank paragraph that separated the image paragraph from the\n#    \"Known Bugs\" heading, effectively moving \"Known Bugs\" directly under it.\n$d.Paragraphs.Item(17).Range.Delete()\n\n# 3. \"Sometimes the keystore will be wiped\" -> \"Sometimes the keystore file will be wiped\"\n$wipeKeystorePara = $d.Paragraphs.Item(20)\n[void]$wipeKeystorePara.Range.Find.Execute(\"keystore will be wiped\", $false, $false, $false, $false, $false, $true, 1, $false, \"keystore file will be wiped\", 2)\n\n# 4. \"save the file\" -> \"save the updated keystore file\"\n$saveFilePara = $d.Paragraphs.Item(21)\n[void]$saveFilePara.Range.Find.Execute(\"save the file\", $false, $false, $false, $false, $false, $true, 1, $false, \"save the updated keystore file\", 2)\n\n# 5. \"temperature the humidity\" -> \"temperature, the humidity\"\n$humidityPara = $d.Paragraphs.Item(25)\n[void]$humidityPara.Range.Find.Execute(\"temperature the humidity\", $false, $false, $false, $false, $false, $true, 1, $false, \"temperature, the humidity\", 2)\n\n# 6. \"so it also outputs\" -> \"so it will output\"\n[void]$humidityPara.Range.Find.Execute(\"so it also outputs\", $false, $false, $false, $false, $false, $true, 1, $false, \"so it will output\", 2)\n"}
